# Autogenerated data refresh for "MSME Country Indicators - Slovenia Summary"
# Updates a handful of percentage/density figures on the "Summary" sheet to
# their more precise values. These cells hold the numbers as text (so the
# trailing zero-less / rounded display is preserved exactly), so we force
# Text formatting before writing each value - otherwise Excel would helpfully
# re-interpret a numeric-looking string like "36.29" as a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 cell -> new value
$updates = @{
    "B11" = "36.29"
    "C11" = "30.67"
    "D11" = "66.96"
    "D12" = "24.71"
    "B33" = "53.75"
    "C33" = "3.37"
    "D33" = "57.12"
    "B34" = "32.93"
    "D34" = "72.92"
    "B36" = "93.92"
    "C36" = "5.88"
    "D36" = "99.81"
    "C40" = "42.36"
    "D40" = "63.56"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
